$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-11-13 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-14 Tuesday", 1) | Out-Null
$d.Content.Find.Execute("70+22=92", $true, $false, $false, $false, $false, $true, 1, $false, "93-23=70", 1) | Out-Null
$d.Content.Find.Execute("85-60=25", $true, $false, $false, $false, $false, $true, 1, $false, "13+25=38", 1) | Out-Null
$d.Content.Find.Execute("29+34=63", $true, $false, $false, $false, $false, $true, 1, $false, "17+73=90", 1) | Out-Null
$d.Content.Find.Execute("88-9=79", $true, $false, $false, $false, $false, $true, 1, $false, "21+11=32", 1) | Out-Null
$d.Content.Find.Execute("75-39=36", $true, $false, $false, $false, $false, $true, 1, $false, "50-50=0", 1) | Out-Null
$d.Content.Find.Execute("57+38=95", $true, $false, $false, $false, $false, $true, 1, $false, "62-16=46", 1) | Out-Null
$d.Content.Find.Execute("3+28=31", $true, $false, $false, $false, $false, $true, 1, $false, "20+14=34", 1) | Out-Null
$d.Content.Find.Execute("31+6=37", $true, $false, $false, $false, $false, $true, 1, $false, "55+15=70", 1) | Out-Null
$d.Content.Find.Execute("21+9=30", $true, $false, $false, $false, $false, $true, 1, $false, "11+55=66", 1) | Out-Null
$d.Content.Find.Execute("28-8=20", $true, $false, $false, $false, $false, $true, 1, $false, "26+50=76", 1) | Out-Null
$d.Content.Find.Execute("31+26=57", $true, $false, $false, $false, $false, $true, 1, $false, "24-14=10", 1) | Out-Null
$d.Content.Find.Execute("7+38=45", $true, $false, $false, $false, $false, $true, 1, $false, "46-43=3", 1) | Out-Null
$d.Content.Find.Execute("77-51=26", $true, $false, $false, $false, $false, $true, 1, $false, "58-4=54", 1) | Out-Null
$d.Content.Find.Execute("14+80=94", $true, $false, $false, $false, $false, $true, 1, $false, "68+3=71", 1) | Out-Null
$d.Content.Find.Execute("41+26=67", $true, $false, $false, $false, $false, $true, 1, $false, "21-2=19", 1) | Out-Null
$d.Content.Find.Execute("14+60=74", $true, $false, $false, $false, $false, $true, 1, $false, "97-64=33", 1) | Out-Null
$d.Content.Find.Execute("95-81=14", $true, $false, $false, $false, $false, $true, 1, $false, "55+21=76", 1) | Out-Null
$d.Content.Find.Execute("34+47=81", $true, $false, $false, $false, $false, $true, 1, $false, "71+16=87", 1) | Out-Null
$d.Content.Find.Execute("93-82=11", $true, $false, $false, $false, $false, $true, 1, $false, "15+56=71", 1) | Out-Null
$d.Content.Find.Execute("14+80=94", $true, $false, $false, $false, $false, $true, 1, $false, "3+83=86", 1) | Out-Null
$d.Content.Find.Execute("75-45=30", $true, $false, $false, $false, $false, $true, 1, $false, "86-56=30", 1) | Out-Null
$d.Content.Find.Execute("85-81=4", $true, $false, $false, $false, $false, $true, 1, $false, "32+12=44", 1) | Out-Null
$d.Content.Find.Execute("37-18=19", $true, $false, $false, $false, $false, $true, 1, $false, "23+20=43", 1) | Out-Null
$d.Content.Find.Execute("74-58=16", $true, $false, $false, $false, $false, $true, 1, $false, "62-18=44", 1) | Out-Null
$d.Content.Find.Execute("51-21=30", $true, $false, $false, $false, $false, $true, 1, $false, "47+44=91", 1) | Out-Null
$d.Content.Find.Execute("59-41=18", $true, $false, $false, $false, $false, $true, 1, $false, "37+60=97", 1) | Out-Null
$d.Content.Find.Execute("73-11=62", $true, $false, $false, $false, $false, $true, 1, $false, "17+2=19", 1) | Out-Null
$d.Content.Find.Execute("20+11=31", $true, $false, $false, $false, $false, $true, 1, $false, "64-49=15", 1) | Out-Null
$d.Content.Find.Execute("63+11=74", $true, $false, $false, $false, $false, $true, 1, $false, "67-18=49", 1) | Out-Null
$d.Content.Find.Execute("27-13=14", $true, $false, $false, $false, $false, $true, 1, $false, "90-71=19", 1) | Out-Null
$d.Content.Find.Execute("67-4=63", $true, $false, $false, $false, $false, $true, 1, $false, "30+31=61", 1) | Out-Null
$d.Content.Find.Execute("46+14=60", $true, $false, $false, $false, $false, $true, 1, $false, "57+30=87", 1) | Out-Null
$d.Content.Find.Execute("54-49=5", $true, $false, $false, $false, $false, $true, 1, $false, "2+59=61", 1) | Out-Null
$d.Content.Find.Execute("52-44=8", $true, $false, $false, $false, $false, $true, 1, $false, "57-23=34", 1) | Out-Null
$d.Content.Find.Execute("72+19=91", $true, $false, $false, $false, $false, $true, 1, $false, "77-75=2", 1) | Out-Null
$d.Content.Find.Execute("46+29=75", $true, $false, $false, $false, $false, $true, 1, $false, "87-32=55", 1) | Out-Null
$d.Content.Find.Execute("35+47=82", $true, $false, $false, $false, $false, $true, 1, $false, "81-53=28", 1) | Out-Null
$d.Content.Find.Execute("91-62=29", $true, $false, $false, $false, $false, $true, 1, $false, "64+29=93", 1) | Out-Null
$d.Content.Find.Execute("20+58=78", $true, $false, $false, $false, $false, $true, 1, $false, "72+2=74", 1) | Out-Null
$d.Content.Find.Execute("10+85=95", $true, $false, $false, $false, $false, $true, 1, $false, "0+84=84", 1) | Out-Null
$d.Content.Find.Execute("33+15=48", $true, $false, $false, $false, $false, $true, 1, $false, "7+55=62", 1) | Out-Null
$d.Content.Find.Execute("1-0=1", $true, $false, $false, $false, $false, $true, 1, $false, "27+46=73", 1) | Out-Null
$d.Content.Find.Execute("37+46=83", $true, $false, $false, $false, $false, $true, 1, $false, "8+17=25", 1) | Out-Null
$d.Content.Find.Execute("38+34=72", $true, $false, $false, $false, $false, $true, 1, $false, "53-6=47", 1) | Out-Null
$d.Content.Find.Execute("26+37=63", $true, $false, $false, $false, $false, $true, 1, $false, "76-32=44", 1) | Out-Null
$d.Content.Find.Execute("56+12=68", $true, $false, $false, $false, $false, $true, 1, $false, "82+4=86", 1) | Out-Null
$d.Content.Find.Execute("35-25=10", $true, $false, $false, $false, $false, $true, 1, $false, "42-31=11", 1) | Out-Null
$d.Content.Find.Execute("98-12=86", $true, $false, $false, $false, $false, $true, 1, $false, "83-41=42", 1) | Out-Null
$d.Content.Find.Execute("35+10=45", $true, $false, $false, $false, $false, $true, 1, $false, "48-8=40", 1) | Out-Null
$d.Content.Find.Execute("69-52=17", $true, $false, $false, $false, $false, $true, 1, $false, "22+63=85", 1) | Out-Null
$d.Content.Find.Execute("32+47=79", $true, $false, $false, $false, $false, $true, 1, $false, "58+20=78", 1) | Out-Null
$d.Content.Find.Execute("54-23=31", $true, $false, $false, $false, $false, $true, 1, $false, "4+87=91", 1) | Out-Null
$d.Content.Find.Execute("28+29=57", $true, $false, $false, $false, $false, $true, 1, $false, "49-40=9", 1) | Out-Null
$d.Content.Find.Execute("39+2=41", $true, $false, $false, $false, $false, $true, 1, $false, "55+25=80", 1) | Out-Null
$d.Content.Find.Execute("83-75=8", $true, $false, $false, $false, $false, $true, 1, $false, "27-8=19", 1) | Out-Null
$d.Content.Find.Execute("96-86=10", $true, $false, $false, $false, $false, $true, 1, $false, "23+69=92", 1) | Out-Null
$d.Content.Find.Execute("22+35=57", $true, $false, $false, $false, $false, $true, 1, $false, "1+71=72", 1) | Out-Null
$d.Content.Find.Execute("67+5=72", $true, $false, $false, $false, $false, $true, 1, $false, "78-62=16", 1) | Out-Null
$d.Content.Find.Execute("17+9=26", $true, $false, $false, $false, $false, $true, 1, $false, "60-46=14", 1) | Out-Null
$d.Content.Find.Execute("84-35=49", $true, $false, $false, $false, $false, $true, 1, $false, "42-9=33", 1) | Out-Null
$d.Content.Find.Execute("46+53=99", $true, $false, $false, $false, $false, $true, 1, $false, "76+6=82", 1) | Out-Null
$d.Content.Find.Execute("90-56=34", $true, $false, $false, $false, $false, $true, 1, $false, "49+8=57", 1) | Out-Null
$d.Content.Find.Execute("39+52=91", $true, $false, $false, $false, $false, $true, 1, $false, "39+8=47", 1) | Out-Null
$d.Content.Find.Execute("62-34=28", $true, $false, $false, $false, $false, $true, 1, $false, "23-2=21", 1) | Out-Null
$d.Content.Find.Execute("62+0=62", $true, $false, $false, $false, $false, $true, 1, $false, "99-87=12", 1) | Out-Null
$d.Content.Find.Execute("40+49=89", $true, $false, $false, $false, $false, $true, 1, $false, "93-48=45", 1) | Out-Null
$d.Content.Find.Execute("64+6=70", $true, $false, $false, $false, $false, $true, 1, $false, "93-12=81", 1) | Out-Null
$d.Content.Find.Execute("74-18=56", $true, $false, $false, $false, $false, $true, 1, $false, "38+4=42", 1) | Out-Null
$d.Content.Find.Execute("20-15=5", $true, $false, $false, $false, $false, $true, 1, $false, "0+95=95", 1) | Out-Null
$d.Content.Find.Execute("84-67=17", $true, $false, $false, $false, $false, $true, 1, $false, "66-59=7", 1) | Out-Null
$d.Content.Find.Execute("72-41=31", $true, $false, $false, $false, $false, $true, 1, $false, "84-14=70", 1) | Out-Null
$d.Content.Find.Execute("82-75=7", $true, $false, $false, $false, $false, $true, 1, $false, "86-7=79", 1) | Out-Null
$d.Content.Find.Execute("75-58=17", $true, $false, $false, $false, $false, $true, 1, $false, "95-43=52", 1) | Out-Null
$d.Content.Find.Execute("64-23=41", $true, $false, $false, $false, $false, $true, 1, $false, "79-37=42", 1) | Out-Null
$d.Content.Find.Execute("9+73=82", $true, $false, $false, $false, $false, $true, 1, $false, "91-79=12", 1) | Out-Null
$d.Content.Find.Execute("59+8=67", $true, $false, $false, $false, $false, $true, 1, $false, "61-51=10", 1) | Out-Null
$d.Content.Find.Execute("26-6=20", $true, $false, $false, $false, $false, $true, 1, $false, "52+8=60", 1) | Out-Null
$d.Content.Find.Execute("3+65=68", $true, $false, $false, $false, $false, $true, 1, $false, "74-7=67", 1) | Out-Null
$d.Content.Find.Execute("71+27=98", $true, $false, $false, $false, $false, $true, 1, $false, "79-38=41", 1) | Out-Null
$d.Content.Find.Execute("21-13=8", $true, $false, $false, $false, $false, $true, 1, $false, "1+23=24", 1) | Out-Null
$d.Content.Find.Execute("37+9=46", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=18", 1) | Out-Null
$d.Content.Find.Execute("16+26=42", $true, $false, $false, $false, $false, $true, 1, $false, "19+68=87", 1) | Out-Null
$d.Content.Find.Execute("58+10=68", $true, $false, $false, $false, $false, $true, 1, $false, "84-62=22", 1) | Out-Null
$d.Content.Find.Execute("51+9=60", $true, $false, $false, $false, $false, $true, 1, $false, "84+10=94", 1) | Out-Null
$d.Content.Find.Execute("5+47=52", $true, $false, $false, $false, $false, $true, 1, $false, "97-70=27", 1) | Out-Null
$d.Content.Find.Execute("50+16=66", $true, $false, $false, $false, $false, $true, 1, $false, "81-18=63", 1) | Out-Null
$d.Content.Find.Execute("97-22=75", $true, $false, $false, $false, $false, $true, 1, $false, "84-82=2", 1) | Out-Null
$d.Content.Find.Execute("46-37=9", $true, $false, $false, $false, $false, $true, 1, $false, "31-28=3", 1) | Out-Null
$d.Content.Find.Execute("70-26=44", $true, $false, $false, $false, $false, $true, 1, $false, "84-51=33", 1) | Out-Null
$d.Content.Find.Execute("71+8=79", $true, $false, $false, $false, $false, $true, 1, $false, "29-14=15", 1) | Out-Null
$d.Content.Find.Execute("68-3=65", $true, $false, $false, $false, $false, $true, 1, $false, "95-64=31", 1) | Out-Null
$d.Content.Find.Execute("15-9=6", $true, $false, $false, $false, $false, $true, 1, $false, "82-63=19", 1) | Out-Null
$d.Content.Find.Execute("31-19=12", $true, $false, $false, $false, $false, $true, 1, $false, "18-0=18", 1) | Out-Null
$d.Content.Find.Execute("79-42=37", $true, $false, $false, $false, $false, $true, 1, $false, "45+32=77", 1) | Out-Null
$d.Content.Find.Execute("93-87=6", $true, $false, $false, $false, $false, $true, 1, $false, "24+16=40", 1) | Out-Null
$d.Content.Find.Execute("64-56=8", $true, $false, $false, $false, $false, $true, 1, $false, "71+9=80", 1) | Out-Null
$d.Content.Find.Execute("83-1=82", $true, $false, $false, $false, $false, $true, 1, $false, "71-7=64", 1) | Out-Null
$d.Content.Find.Execute("10+58=68", $true, $false, $false, $false, $false, $true, 1, $false, "97-42=55", 1) | Out-Null
$d.Content.Find.Execute("55+26=81", $true, $false, $false, $false, $false, $true, 1, $false, "0+40=40", 1) | Out-Null
$d.Content.Find.Execute("74-74=0", $true, $false, $false, $false, $false, $true, 1, $false, "72-4=68", 1) | Out-Null
